$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to a number by Excel;
# force a Text number format, assign the value, then restore the default "Normal"
# style so the cell keeps no explicit style (matching the original formatting).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "240.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.927"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.573"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "69.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "40.98"
$ws.Range("D51").Style = "Normal"

# Remaining cells: plain text assignment (Excel keeps these as text naturally).
$ws.Range("D2").Value = "27.596.08"
$ws.Range("E2").Value = "  -2.54%  "
$ws.Range("D3").Value = "1.657.35"
$ws.Range("E3").Value = "  -4.33%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -2.74%  "
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("D12").Value = "1.892.41"
$ws.Range("E12").Value = "  -4.32%  "
$ws.Range("D13").Value = "1.668.13"
$ws.Range("E13").Value = "  -3.63%  "
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "27.560.43"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("E20").Value = "  -4.40%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -4.14%  "
$ws.Range("E23").Value = "  -4.29%  "
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("E26").Value = "  -4.67%  "
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").Value = "1.459.88"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("E34").Value = "  -5.16%  "
$ws.Range("E35").Value = "  -5.25%  "
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("E37").Value = "  -5.41%  "
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("E39").Value = "  -5.29%  "
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("E41").Value = "  -3.97%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -4.15%  "
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").Value = "1.800.32"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("E49").Value = "  -6.13%  "
$ws.Range("E50").Value = "  -1.85%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("E51").Value = "  +11.84%  "
